# Updates the cached display text of the "slide number" field placeholder
# (‹nº›  ->  ‹#›) on every slide layout, the slide master, and the notes
# master of the active presentation. The nine real slides already carry a
# live, already-computed slide number (e.g. "2", "3", ...) in their field
# caches and are intentionally left untouched.

$p = $ppt.ActivePresentation

$oldText = [string][char]0x2039 + "n" + [char]0xBA + [char]0x203A   # ‹nº›
$newText = [string][char]0x2039 + "#" + [char]0x203A                # ‹#›

# ppPlaceholderSlideNumber
$SlideNumberPlaceholder = 13

function Update-SlideNumberField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)

        $phType = -1
        try { $phType = $shape.PlaceholderFormat.Type } catch {}

        if ($phType -eq $SlideNumberPlaceholder -and $shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}

# 1) Every slide layout ("CustomLayout") hanging off the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-SlideNumberField $layouts.Item($li).Shapes
}

# 2) The slide master itself.
Update-SlideNumberField $p.SlideMaster.Shapes

# 3) The notes master.
Update-SlideNumberField $p.NotesMaster.Shapes
